#
# Stock reconciliation corrections for CryCompanywiseStockReport_1.xlsx
# For each affected item row: update Qty (F) and/or recompute Value (G = Rate(D) * Qty(F)).
# A few adjacent item-row pairs were reordered/swapped (e.g. rows 277/278, 476/477,
# 487/488, 710/711, 737/738) -- their Code/Rate/Qty/Value cells are updated to reflect
# the corrected row order. Each "Sub Total:" / "Grand Total:" row (column B) is updated
# to match the corrected sum of its section.
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F95: 8->7, G95: 2961.44->2591.26
$ws.Cells.Item(95, 6).Value = 7
$ws.Cells.Item(95, 7).Value = 2591.26

# F104: 3->2, G104: 307.38->204.92
$ws.Cells.Item(104, 6).Value = 2
$ws.Cells.Item(104, 7).Value = 204.92

# B114: 242493.4->242020.76
$ws.Cells.Item(114, 2).Value = 242020.76

# F193: 289->288, G193: 18727.2->18662.4
$ws.Cells.Item(193, 6).Value = 288
$ws.Cells.Item(193, 7).Value = 18662.4

# B200: 43725.58->43660.78
$ws.Cells.Item(200, 2).Value = 43660.78

# F265: 22->21, G265: 3685->3517.5
$ws.Cells.Item(265, 6).Value = 21
$ws.Cells.Item(265, 7).Value = 3517.5

# B274: 71413.22->71245.72
$ws.Cells.Item(274, 2).Value = 71245.72

# B277: 61610->63565, E277: 122.71->109.19, F277: -58->60, G277: -5957.18->6162.6
$ws.Cells.Item(277, 2).Value = 63565
$ws.Cells.Item(277, 5).Value = 109.19
$ws.Cells.Item(277, 6).Value = 60
$ws.Cells.Item(277, 7).Value = 6162.6

# B278: 63565->61610, E278: 109.19->122.71, F278: 60->-58, G278: 6162.6->-5957.18
$ws.Cells.Item(278, 2).Value = 61610
$ws.Cells.Item(278, 5).Value = 122.71
$ws.Cells.Item(278, 6).Value = -58
$ws.Cells.Item(278, 7).Value = -5957.18

# F287: 22->24, G287: 2914.56->3179.52
$ws.Cells.Item(287, 6).Value = 24
$ws.Cells.Item(287, 7).Value = 3179.52

# F288: 3->1, G288: 868.83->289.61
$ws.Cells.Item(288, 6).Value = 1
$ws.Cells.Item(288, 7).Value = 289.61

# F328: 381->355, G328: 8012.43->7465.65
$ws.Cells.Item(328, 6).Value = 355
$ws.Cells.Item(328, 7).Value = 7465.65

# F330: 5->4, G330: 2629.75->2103.8
$ws.Cells.Item(330, 6).Value = 4
$ws.Cells.Item(330, 7).Value = 2103.8

# B339: 274339.13->272952.14
$ws.Cells.Item(339, 2).Value = 272952.14

# F363: 24->20, G363: 510->425
$ws.Cells.Item(363, 6).Value = 20
$ws.Cells.Item(363, 7).Value = 425

# B395: 232630.46->232545.46
$ws.Cells.Item(395, 2).Value = 232545.46

# F426: 53->52, G426: 5119.8->5023.2
$ws.Cells.Item(426, 6).Value = 52
$ws.Cells.Item(426, 7).Value = 5023.2

# B430: 40881.29->40784.69
$ws.Cells.Item(430, 2).Value = 40784.69

# B476: 64922->45706, E476: 20.98->23.58, F476: 68->-202, G476: 1341.64->-3985.46
$ws.Cells.Item(476, 2).Value = 45706
$ws.Cells.Item(476, 5).Value = 23.58
$ws.Cells.Item(476, 6).Value = -202
$ws.Cells.Item(476, 7).Value = -3985.46

# B477: 45706->64922, E477: 23.58->20.98, F477: -202->68, G477: -3985.46->1341.64
$ws.Cells.Item(477, 2).Value = 64922
$ws.Cells.Item(477, 5).Value = 20.98
$ws.Cells.Item(477, 6).Value = 68
$ws.Cells.Item(477, 7).Value = 1341.64

# F484: 434->431, G484: 2816.66->2797.19
$ws.Cells.Item(484, 6).Value = 431
$ws.Cells.Item(484, 7).Value = 2797.19

# F486: 149->147, G486: 1959.35->1933.05
$ws.Cells.Item(486, 6).Value = 147
$ws.Cells.Item(486, 7).Value = 1933.05

# B487: 45702->64919, E487: 31.43->27.97, F487: -215->65, G487: -5654.5->1709.5
$ws.Cells.Item(487, 2).Value = 64919
$ws.Cells.Item(487, 5).Value = 27.97
$ws.Cells.Item(487, 6).Value = 65
$ws.Cells.Item(487, 7).Value = 1709.5

# B488: 64919->45702, E488: 27.97->31.43, F488: 66->-215, G488: 1735.8->-5654.5
$ws.Cells.Item(488, 2).Value = 45702
$ws.Cells.Item(488, 5).Value = 31.43
$ws.Cells.Item(488, 6).Value = -215
$ws.Cells.Item(488, 7).Value = -5654.5

# F489: 30->27, G489: 492.9->443.61
$ws.Cells.Item(489, 6).Value = 27
$ws.Cells.Item(489, 7).Value = 443.61

# F490: 212->210, G490: 3122.76->3093.3
$ws.Cells.Item(490, 6).Value = 210
$ws.Cells.Item(490, 7).Value = 3093.3

# B492: -12355.38->-12506.2
$ws.Cells.Item(492, 2).Value = -12506.2

# F497: 2->0, G497: 100.16->0
$ws.Cells.Item(497, 6).Value = 0
$ws.Cells.Item(497, 7).Value = 0

# B508: 7990.3->7890.14
$ws.Cells.Item(508, 2).Value = 7890.14

# F516: 10->9, G516: 161.8->145.62
$ws.Cells.Item(516, 6).Value = 9
$ws.Cells.Item(516, 7).Value = 145.62

# B528: 16791.72->16775.54
$ws.Cells.Item(528, 2).Value = 16775.54

# F545: 25->24, G545: 6136.75->5891.28
$ws.Cells.Item(545, 6).Value = 24
$ws.Cells.Item(545, 7).Value = 5891.28

# F546: 2->1, G546: 5558.8->2779.4
$ws.Cells.Item(546, 6).Value = 1
$ws.Cells.Item(546, 7).Value = 2779.4

# B547: 16931.12->13906.25
$ws.Cells.Item(547, 2).Value = 13906.25

# F549: 280->278, G549: 1904->1890.4
$ws.Cells.Item(549, 6).Value = 278
$ws.Cells.Item(549, 7).Value = 1890.4

# F551: 105->104, G551: 1353.45->1340.56
$ws.Cells.Item(551, 6).Value = 104
$ws.Cells.Item(551, 7).Value = 1340.56

# B557: 6304.49->6278
$ws.Cells.Item(557, 2).Value = 6278

# F616: 52->51, G616: 9124.44->8948.97
$ws.Cells.Item(616, 6).Value = 51
$ws.Cells.Item(616, 7).Value = 8948.97

# B619: 38067.77->37892.3
$ws.Cells.Item(619, 2).Value = 37892.3

# F655: 317->315, G655: 25480.46->25319.7
$ws.Cells.Item(655, 6).Value = 315
$ws.Cells.Item(655, 7).Value = 25319.7

# B656: 33881.22->33720.46
$ws.Cells.Item(656, 2).Value = 33720.46

# F706: 117->116, G706: 16746.21->16603.08
$ws.Cells.Item(706, 6).Value = 116
$ws.Cells.Item(706, 7).Value = 16603.08

# B710: 61428->63150, D710: 69.16->75.68, E710: 73.52->80.45, F710: 1->33, G710: 69.16->2497.44
$ws.Cells.Item(710, 2).Value = 63150
$ws.Cells.Item(710, 4).Value = 75.68
$ws.Cells.Item(710, 5).Value = 80.45
$ws.Cells.Item(710, 6).Value = 33
$ws.Cells.Item(710, 7).Value = 2497.44

# B711: 63150->61428, D711: 75.68->69.16, E711: 80.45->73.52, F711: 33->1, G711: 2497.44->69.16
$ws.Cells.Item(711, 2).Value = 61428
$ws.Cells.Item(711, 4).Value = 69.16
$ws.Cells.Item(711, 5).Value = 73.52
$ws.Cells.Item(711, 6).Value = 1
$ws.Cells.Item(711, 7).Value = 69.16

# F715: 75->74, G715: 5217->5147.44
$ws.Cells.Item(715, 6).Value = 74
$ws.Cells.Item(715, 7).Value = 5147.44

# F720: 164->160, G720: 19796.44->19313.6
$ws.Cells.Item(720, 6).Value = 160
$ws.Cells.Item(720, 7).Value = 19313.6

# B721: 90944.18->90248.65
$ws.Cells.Item(721, 2).Value = 90248.65

# B737: 65362->65079, F737: 44->21, G737: 1798.28->858.27
$ws.Cells.Item(737, 2).Value = 65079
$ws.Cells.Item(737, 6).Value = 21
$ws.Cells.Item(737, 7).Value = 858.27

# B738: 65079->65362, F738: 21->44, G738: 858.27->1798.28
$ws.Cells.Item(738, 2).Value = 65362
$ws.Cells.Item(738, 6).Value = 44
$ws.Cells.Item(738, 7).Value = 1798.28

# F740: 62->60, G740: 2049.72->1983.6
$ws.Cells.Item(740, 6).Value = 60
$ws.Cells.Item(740, 7).Value = 1983.6

# F742: 134->131, G742: 5760.66->5631.69
$ws.Cells.Item(742, 6).Value = 131
$ws.Cells.Item(742, 7).Value = 5631.69

# F746: 90->89, G746: 21762->21520.2
$ws.Cells.Item(746, 6).Value = 89
$ws.Cells.Item(746, 7).Value = 21520.2

# F747: 100->99, G747: 5699->5642.01
$ws.Cells.Item(747, 6).Value = 99
$ws.Cells.Item(747, 7).Value = 5642.01

# B748: 64020.34->63526.46
$ws.Cells.Item(748, 2).Value = 63526.46

# F775: 558->556, G775: 157841.46->157275.72
$ws.Cells.Item(775, 6).Value = 556
$ws.Cells.Item(775, 7).Value = 157275.72

# F779: 90->88, G779: 11572.2->11315.04
$ws.Cells.Item(779, 6).Value = 88
$ws.Cells.Item(779, 7).Value = 11315.04

# B780: 703261.3->702438.4
$ws.Cells.Item(780, 2).Value = 702438.4

# B798: 2649692.95->2641752.36
$ws.Cells.Item(798, 2).Value = 2641752.36

# B799: 2649692.95->2641752.36
$ws.Cells.Item(799, 2).Value = 2641752.36
